# Auto Appium Server & Emulator Detector Enhancement
# - Fix the TC003 "Scroll down..." title text on the "TCs" sheet (drop the
#   trailing period).
# - Refresh the last-used selection/active-cell on each sheet.

$wb = $excel.ActiveWorkbook

$wsTCs = $wb.Worksheets.Item("TCs")
$wsRadio = $wb.Worksheets.Item("RadioButtons")

# Correct the trailing period on the "Scroll down..." test-case title.
$wsTCs.Range("B4").Value = "Scroll down to the bottom and uncheck the show notification checkbox"

# Restore the last-recorded selections for each sheet.
$wsTCs.Activate()
$wsTCs.Range("B8").Select() | Out-Null

$wsRadio.Activate()
$wsRadio.Range("B12").Select() | Out-Null
